$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (Number of appeals to the Ombudsman, 2007-2020) gains three more
# year columns (2021, 2022, 2023), shifting the used range from Q to T.
# Copy the formatting of the existing last data column (Q, rows 2-5) into the
# new R:T columns so borders/fonts/number formats match, then fill in the
# new year headers and data values.
$ws.Range("Q2:Q5").Copy()
$ws.Range("R2:T2").PasteSpecial(-4122)

$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023

$ws.Range("R4").Value = 4301
$ws.Range("S4").Value = 3690
$ws.Range("T4").Value = 2620

$ws.Range("R5").Value = 427
$ws.Range("S5").Value = 280
$ws.Range("T5").Value = 264

# Reset the selection to the sheet's top-left cell (clears the stray
# F16 selection that was saved in the original file).
$ws.Range("A1").Select()
